$d = $word.ActiveDocument
$d.Content.Find.Execute("Password: cakePIE480", $true, $false, $false, $false, $false, $true, 1, $false, "Password:480cakePIE", 2)

# split into two runs first
$p2 = $d.Paragraphs(2)
$splitPos = $p2.Range.Start + 9
$tailRange = $d.Range($splitPos, $p2.Range.End - 1)
$d.Bookmarks.Add("__tmp_split", $tailRange)
$d.Bookmarks("__tmp_split").Delete()

Write-Host "p2 text:" $d.Paragraphs(2).Range.Text
Write-Host "p2 start/end:" $d.Paragraphs(2).Range.Start $d.Paragraphs(2).Range.End

# Now add _GoBack spanning exactly the second run "480cakePIE" (splitPos to End-1)
$p2b = $d.Paragraphs(2)
$fullRange = $d.Range($splitPos, $p2b.Range.End - 1)
Write-Host "fullRange text:[" $fullRange.Text "]"
$d.Bookmarks.Add("_GoBack", $fullRange)
